$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes (XML "width" is ColumnWidth + 5/6; subtract that
# offset so the saved width attribute matches the target values exactly).
$ws.Columns.Item(1).ColumnWidth = 83 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 35 - (5/6)

# Text content changes
$ws.Range("B1").Value = "div_testRuns_internalRoleCellName"
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleTestRunWithEnvironmentAndDevices-test-data"
